$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.992518924015428
$ws.Range("C2").Value = 0.2843258085016487
$ws.Range("D2").Value = 0.05989484274526546
$ws.Range("F2").Value = 1.703210617878142
$ws.Range("G2").Value = 0.002481720555312002
$ws.Range("I2").Value = 1.222840110520707
$ws.Range("L2").Value = 0.303207780257253
$ws.Range("N2").Value = 1.482901876590123
$ws.Range("B3").Value = 1.85173105524774
$ws.Range("C3").Value = 0.2474039262088752
$ws.Range("D3").Value = 0.06064901790708532
$ws.Range("F3").Value = 1.665365419445251
$ws.Range("G3").Value = 0.002487219738402047
$ws.Range("I3").Value = 1.213911196093122
$ws.Range("L3").Value = 0.2921938689453896
$ws.Range("N3").Value = 1.501981537002482
$ws.Range("B4").Value = 1.766342198679638
$ws.Range("C4").Value = 0.224743570430519
$ws.Range("D4").Value = 0.06114221135523046
$ws.Range("F4").Value = 1.643352727609098
$ws.Range("G4").Value = 0.002490773090476534
$ws.Range("I4").Value = 1.209234862149444
$ws.Range("L4").Value = 0.2856215616866535
$ws.Range("N4").Value = 1.514309201118099
$ws.Range("B5").Value = 1.731809271615646
$ws.Range("C5").Value = 0.2155110013733577
$ws.Range("D5").Value = 0.0613507361403105
$ws.Range("F5").Value = 1.63468826886627
$ws.Range("G5").Value = 0.002492265731465926
$ws.Range("I5").Value = 1.207530896178753
$ws.Range("L5").Value = 0.2829909156341586
$ws.Range("N5").Value = 1.519486457241012
$ws.Range("B6").Value = 1.726090985595306
$ws.Range("C6").Value = 0.2139780215191536
$ws.Range("D6").Value = 0.06138581628371576
$ws.Range("F6").Value = 1.633267959228291
$ws.Range("G6").Value = 0.002492516283301176
$ws.Range("I6").Value = 1.207260105409681
$ws.Range("L6").Value = 0.2825569686676204
$ws.Range("N6").Value = 1.520355403262016
$ws.Range("B7").Value = 1.7658754099362
$ws.Range("C7").Value = 0.2246190508435575
$ws.Range("D7").Value = 0.06114499308766597
$ws.Range("F7").Value = 1.643234639792212
$ws.Range("G7").Value = 0.002490793039834181
$ws.Range("I7").Value = 1.209211066630779
$ws.Range("L7").Value = 0.2855858913435725
$ws.Range("N7").Value = 1.514378401936817
$ws.Range("B8").Value = 1.943755136022844
$ws.Range("C8").Value = 0.2715926321211271
$ws.Range("D8").Value = 0.06014861069445843
$ws.Range("F8").Value = 1.689906243380435
$ws.Range("G8").Value = 0.002483580066334891
$ws.Range("I8").Value = 1.21959350665329
$ws.Range("L8").Value = 0.2993705508612692
$ws.Range("N8").Value = 1.489353125300713
$ws.Range("B9").Value = 2.301048416411447
$ws.Range("C9").Value = 0.3638246191948724
$ws.Range("D9").Value = 0.05843514426016583
$ws.Range("F9").Value = 1.791241835302444
$ws.Range("G9").Value = 0.00247083140511295
$ws.Range("I9").Value = 1.246397714740482
$ws.Range("L9").Value = 0.3279240392201217
$ws.Range("N9").Value = 1.445158510515292
$ws.Range("B10").Value = 2.568881467237702
$ws.Range("C10").Value = 0.431719295250673
$ws.Range("D10").Value = 0.05732474768618445
$ws.Range("F10").Value = 1.871821087959773
$ws.Range("G10").Value = 0.002462305994638231
$ws.Range("I10").Value = 1.270090948852413
$ws.Range("L10").Value = 0.3498494415566569
$ws.Range("N10").Value = 1.415688208138867
$ws.Range("B11").Value = 2.691921100182128
$ws.Range("C11").Value = 0.4626485606024744
$ws.Range("D11").Value = 0.05685226367103269
$ws.Range("F11").Value = 1.909841806573098
$ws.Range("G11").Value = 0.002458608052386207
$ws.Range("I11").Value = 1.281754116621698
$ws.Range("L11").Value = 0.3600339627088402
$ws.Range("N11").Value = 1.402938094530562
$ws.Range("B12").Value = 2.73868825859654
$ws.Range("C12").Value = 0.4743679225598498
$ws.Range("D12").Value = 0.05667807693744109
$ws.Range("F12").Value = 1.924437998398105
$ws.Range("G12").Value = 0.002457233502821022
$ws.Range("I12").Value = 1.286299156235586
$ws.Range("L12").Value = 0.3639211638809314
$ws.Range("N12").Value = 1.398204798058202
$ws.Range("B13").Value = 2.728608310921118
$ws.Range("C13").Value = 0.4718436179906007
$ws.Range("D13").Value = 0.05671538004924059
$ws.Range("F13").Value = 1.921285578666868
$ws.Range("G13").Value = 0.002457528392215766
$ws.Range("I13").Value = 1.285314568613458
$ws.Range("L13").Value = 0.3630826227809507
$ws.Range("N13").Value = 1.399219968772115
$ws.Range("B14").Value = 2.695765149085389
$ws.Range("C14").Value = 0.4636125719678148
$ws.Range("D14").Value = 0.05683783814690102
$ws.Range("F14").Value = 1.911038651044493
$ws.Range("G14").Value = 0.002458494451546115
$ws.Range("I14").Value = 1.282125458543163
$ws.Range("L14").Value = 0.3603531513863487
$ws.Range("N14").Value = 1.402546778352686
$ws.Range("B15").Value = 2.67567057846145
$ws.Range("C15").Value = 0.4585717712248538
$ws.Range("D15").Value = 0.05691346477505377
$ws.Range("F15").Value = 1.904788040964007
$ws.Range("G15").Value = 0.002459089543681253
$ws.Range("I15").Value = 1.280188800250599
$ws.Range("L15").Value = 0.3586852571249324
$ws.Range("N15").Value = 1.404596919932622
$ws.Range("B16").Value = 2.560864876533913
$ws.Range("C16").Value = 0.4296989450795081
$ws.Range("D16").Value = 0.05735628576199758
$ws.Range("F16").Value = 1.869364005046663
$ws.Range("G16").Value = 0.002462551279506398
$ws.Range("I16").Value = 1.269346643129737
$ws.Range("L16").Value = 0.3491881186360928
$ws.Range("N16").Value = 1.416534706927862
$ws.Range("B17").Value = 2.490744325430228
$ws.Range("C17").Value = 0.4119981752739363
$ws.Range("D17").Value = 0.05763632972944777
$ws.Range("F17").Value = 1.847983702904315
$ws.Range("G17").Value = 0.002464721018958822
$ws.Range("I17").Value = 1.262922841127462
$ws.Range("L17").Value = 0.3434160464378948
$ws.Range("N17").Value = 1.424026578883947
$ws.Range("B18").Value = 2.450525762471557
$ws.Range("C18").Value = 0.4018212165735804
$ws.Range("D18").Value = 0.05780047355884221
$ws.Range("F18").Value = 1.835814641036592
$ws.Range("G18").Value = 0.002465985975955076
$ws.Range("I18").Value = 1.259311268107822
$ws.Range("L18").Value = 0.3401159043852005
$ws.Range("N18").Value = 1.42839742477365
$ws.Range("B19").Value = 2.436927782378234
$ws.Range("C19").Value = 0.3983761442903528
$ws.Range("D19").Value = 0.05785657592907256
$ws.Range("F19").Value = 1.831716375326408
$ws.Range("G19").Value = 0.002466417189638861
$ws.Range("I19").Value = 1.258102710029192
$ws.Range("L19").Value = 0.3390019257954009
$ws.Range("N19").Value = 1.429887901755976
$ws.Range("B20").Value = 2.498197076252666
$ws.Range("C20").Value = 0.4138820265814047
$ws.Range("D20").Value = 0.05760620056286569
$ws.Range("F20").Value = 1.850246375035482
$ws.Range("G20").Value = 0.002464488290187892
$ws.Range("I20").Value = 1.263598043472939
$ws.Range("L20").Value = 0.3440284421048148
$ws.Range("N20").Value = 1.423222664645564
$ws.Range("B21").Value = 2.705407223693385
$ws.Range("C21").Value = 0.4660300296580431
$ws.Range("D21").Value = 0.05680174047628483
$ws.Range("F21").Value = 1.91404301433073
$ws.Range("G21").Value = 0.00245820999815311
$ws.Range("I21").Value = 1.28305868178677
$ws.Range("L21").Value = 0.3611540321805791
$ws.Range("N21").Value = 1.401567033869028
$ws.Range("B22").Value = 2.841850406931485
$ws.Range("C22").Value = 0.5001536695018558
$ws.Range("D22").Value = 0.0563035837617214
$ws.Range("F22").Value = 1.956896325101155
$ws.Range("G22").Value = 0.002454256974226179
$ws.Range("I22").Value = 1.29652660904344
$ws.Range("L22").Value = 0.372524757584614
$ws.Range("N22").Value = 1.387967209947895
$ws.Range("B23").Value = 2.768934160204878
$ws.Range("C23").Value = 0.4819371235270182
$ws.Range("D23").Value = 0.0565669204958219
$ws.Range("F23").Value = 1.933917929597698
$ws.Range("G23").Value = 0.002456353082363958
$ws.Range("I23").Value = 1.289269566712818
$ws.Range("L23").Value = 0.3664395956533326
$ws.Range("N23").Value = 1.395174872479814
$ws.Range("B24").Value = 2.494827392862248
$ws.Range("C24").Value = 0.4130303391831944
$ws.Range("D24").Value = 0.0576198121725966
$ws.Range("F24").Value = 1.849223038773687
$ws.Range("G24").Value = 0.00246459345197189
$ws.Range("I24").Value = 1.263292530552562
$ws.Range("L24").Value = 0.343751521046002
$ws.Range("N24").Value = 1.423585916000008
$ws.Range("B25").Value = 2.203467285830527
$ws.Range("C25").Value = 0.3388546898735285
$ws.Range("D25").Value = 0.05887274697538913
$ws.Range("F25").Value = 1.762762380218931
$ws.Range("G25").Value = 0.002474131843133928
$ws.Range("I25").Value = 1.238449218272478
$ws.Range("L25").Value = 0.3200346550794251
$ws.Range("N25").Value = 1.45658939354789
